# Apply the automatic daily update to the EPEX spot / Gaz / CO2 workbook.
#
# Sheet "Prix Spot": add a new date column (CR) for "17-sep" with the 24
#   hourly prices.
# Sheet "Gaz": append a new row (93) for date 2025-09-15.
# Sheet "CO2": append a new row (93) for date 2025-09-15.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" -> new column CR ("17-sep")
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting (bold / centered / bordered) from the previous
# date column (CQ) onto the new one (CR), then overwrite its text.
$wsPrix.Range("CQ1").Copy($wsPrix.Range("CR1"))
$wsPrix.Range("CR1").Value = "17-sep"

$prixValues = @(
    71.56999999999999,
    63.37,
    50.98,
    40.32,
    32,
    35.16,
    67.78,
    111.07,
    112.08,
    85.84999999999999,
    57.9,
    37.6,
    18.05,
    7.88,
    0.09,
    6.11,
    12.63,
    50,
    83.12,
    111.76,
    123.55,
    95.2,
    87.36,
    63.02
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 96).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" -> new row 93 (2025-09-15)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the new date cell to stay plain text (matching every other date
# cell in the column) instead of being auto-converted into a date serial.
$wsGaz.Range("A93").NumberFormat = "@"
$wsGaz.Range("A93").Value = "2025-09-15"
$wsGaz.Range("A93").Style = "Normal"

$wsGaz.Range("B93").Value = 31.55

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" -> new row 93 (2025-09-15)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A93").NumberFormat = "@"
$wsCo2.Range("A93").Value = "2025-09-15"
$wsCo2.Range("A93").Style = "Normal"

$wsCo2.Range("B93").Value = 76.23999999999999
